$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5095059
$ws.Range("I74").Value = 5660743
$ws.Range("J74").Value = 3900
$ws.Range("K74").Value = 5660743
$ws.Range("L74").Value = 3900
$ws.Range("M74").Value = -5659807
$ws.Range("N74").Value = -5772

$ws.Range("H77").Value = 5095059
$ws.Range("I77").Value = 5660743
$ws.Range("J77").Value = 3900
$ws.Range("K77").Value = 28303715
$ws.Range("L77").Value = 19500
$ws.Range("M77").Value = -28299035
$ws.Range("N77").Value = -28860

$ws.Range("H108").Value = 59684
$ws.Range("J108").Value = 59684
$ws.Range("L108").Value = 59684
$ws.Range("N108").Value = -67364

$ws.Range("H129").Value = 1684696.5
$ws.Range("I129").Value = 548.5
$ws.Range("J129").Value = 1853111.2
$ws.Range("K129").Value = 1645.5
$ws.Range("L129").Value = 5559333.6
$ws.Range("M129").Value = 3354.5
$ws.Range("N129").Value = -5569333.6

$ws.Range("H132").Value = 1772.8429
$ws.Range("I132").Value = 1641.2063
$ws.Range("J132").Value = 2957.5715
$ws.Range("K132").Value = 4923.6189
$ws.Range("L132").Value = 8872.7145
$ws.Range("M132").Value = -2393.6189
$ws.Range("N132").Value = -13932.7145

$ws.Range("H135").Value = 2349.75
$ws.Range("I135").Value = 1199.7273
$ws.Range("J135").Value = 15000
$ws.Range("K135").Value = 10797.5457
$ws.Range("L135").Value = 135000
$ws.Range("M135").Value = -8262.545700000001
$ws.Range("N135").Value = -140070

$ws.Range("H136").Value = 37514.234
$ws.Range("J136").Value = 36716.133
$ws.Range("L136").Value = 36716.133
$ws.Range("N136").Value = -46916.133

$ws.Range("H137").Value = 907.89655
$ws.Range("I137").Value = 852.8570999999999
$ws.Range("J137").Value = 1052.375
$ws.Range("K137").Value = 2558.5713
$ws.Range("L137").Value = 3157.125
$ws.Range("M137").Value = -8.57129999999961
$ws.Range("N137").Value = -8257.125

$ws.Range("H138").Value = 3911.8
$ws.Range("I138").Value = 832.619
$ws.Range("J138").Value = 4730.3164
$ws.Range("K138").Value = 2497.857
$ws.Range("L138").Value = 14190.9492
$ws.Range("M138").Value = 2642.143
$ws.Range("N138").Value = -24470.9492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5654854
$ws.Range("I32").Value = 6414969
$ws.Range("K32").Value = 6414969
$ws.Range("M32").Value = -6414682

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2060.5483
$ws.Range("I86").Value = 1824.7826
$ws.Range("J86").Value = 2738.375
$ws.Range("K86").Value = 1824.7826
$ws.Range("L86").Value = 2738.375
$ws.Range("M86").Value = -701.7826
$ws.Range("N86").Value = -4984.375

$ws.Range("H89").Value = 2060.5483
$ws.Range("I89").Value = 1824.7826
$ws.Range("J89").Value = 2738.375
$ws.Range("K89").Value = 9123.913
$ws.Range("L89").Value = 13691.875
$ws.Range("M89").Value = -3507.913
$ws.Range("N89").Value = -24923.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2040.2653
$ws.Range("I31").Value = 2226.389
$ws.Range("J31").Value = 1524.8462
$ws.Range("K31").Value = 2226.389
$ws.Range("L31").Value = 1524.8462
$ws.Range("M31").Value = -1931.389
$ws.Range("N31").Value = -2114.8462

$ws.Range("H34").Value = 2040.2653
$ws.Range("I34").Value = 2226.389
$ws.Range("J34").Value = 1524.8462
$ws.Range("K34").Value = 2226.389
$ws.Range("L34").Value = 1524.8462
$ws.Range("M34").Value = -2024.389
$ws.Range("N34").Value = -1928.8462

$ws.Range("H62").Value = 2763.8
$ws.Range("I62").Value = 2926.25
$ws.Range("K62").Value = 2926.25
$ws.Range("M62").Value = -2302.25

$ws.Range("H65").Value = 2763.8
$ws.Range("I65").Value = 2926.25
$ws.Range("K65").Value = 14631.25
$ws.Range("M65").Value = -11511.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1066.6666
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4774.75
$ws.Range("I70").Value = 4042.5715
$ws.Range("K70").Value = 4042.5715
$ws.Range("M70").Value = -3772.5715

$ws.Range("H73").Value = 4774.75
$ws.Range("I73").Value = 4042.5715
$ws.Range("K73").Value = 4042.5715
$ws.Range("M73").Value = -3106.5715

$ws.Range("H133").Value = 38653.223
$ws.Range("J133").Value = 38653.223
$ws.Range("L133").Value = 38653.223
$ws.Range("N133").Value = -48773.223

$ws.Range("H135").Value = 12760
$ws.Range("J135").Value = 12760
$ws.Range("L135").Value = 12760
$ws.Range("N135").Value = -22900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1187.1333
$ws.Range("J46").Value = 965.4
$ws.Range("L46").Value = 965.4
$ws.Range("N46").Value = -1341.4

$ws.Range("H132").Value = 1989.3478
$ws.Range("I132").Value = 1343.4667
$ws.Range("K132").Value = 4030.4001
$ws.Range("M132").Value = -1500.4001

$ws.Range("H136").Value = 3941.2104
$ws.Range("I136").Value = 1875.2307
$ws.Range("K136").Value = 5625.6921
$ws.Range("M136").Value = -3075.6921

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 44262.453
$ws.Range("J46").Value = 44262.453
$ws.Range("L46").Value = 44262.453
$ws.Range("N46").Value = -44724.453

$ws.Range("H93").Value = 29889
$ws.Range("J93").Value = 29889
$ws.Range("L93").Value = 29889
$ws.Range("N93").Value = -34881

$ws.Range("H101").Value = 20301
$ws.Range("J101").Value = 20301
$ws.Range("L101").Value = 20301
$ws.Range("N101").Value = -26791

$ws.Range("H105").Value = 50000
$ws.Range("J105").Value = 50000
$ws.Range("L105").Value = 50000
$ws.Range("N105").Value = -56988

$ws.Range("H107").Value = 569.4828
$ws.Range("I107").Value = 564.5
$ws.Range("J107").Value = 577.63635
$ws.Range("K107").Value = 1693.5
$ws.Range("L107").Value = 1732.90905
$ws.Range("M107").Value = 226.5
$ws.Range("N107").Value = -5572.90905

$ws.Range("H113").Value = 327.55554
$ws.Range("J113").Value = 350.4
$ws.Range("L113").Value = 1051.2
$ws.Range("N113").Value = -5391.2

$ws.Range("H134").Value = 44262.453
$ws.Range("J134").Value = 44262.453
$ws.Range("L134").Value = 132787.359
$ws.Range("N134").Value = -137857.359

$ws.Range("H135").Value = 78376.5
$ws.Range("J135").Value = 78376.5
$ws.Range("L135").Value = 78376.5
$ws.Range("N135").Value = -88516.5

$ws.Range("H136").Value = 1483.2632
$ws.Range("I136").Value = 1557.2683
$ws.Range("K136").Value = 4671.8049
$ws.Range("M136").Value = -2121.8049
